$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136, pushing existing rows 136-232 down to 137-233.
$ws.Rows.Item(136).Insert()

# Populate the newly inserted row 136 with the new weekly data point.
$ws.Cells.Item(136, 1).Value = 10
$ws.Cells.Item(136, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(136, 3).Value = "La Araucanía"
$ws.Cells.Item(136, 4).Value = 44957
$ws.Cells.Item(136, 5).Value = 9
$ws.Cells.Item(136, 6).Value = 100112012
$ws.Cells.Item(136, 7).Value = "Espinaca"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 40
$ws.Cells.Item(136, 11).Value = 12000
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = 12000
$ws.Cells.Item(136, 14).Value = "$/docena de atados"
$ws.Cells.Item(136, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(136, 16).Value = 4000
$ws.Cells.Item(136, 17).Value = 3
$ws.Cells.Item(136, 18).Value = "Hortaliza"
